$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A (a numeric "count" column with a bordered/centered style)
# is removed entirely. Deleting the whole column shifts B:F left by one, so
# the former B:F data (headers QS_Astral5/FNRATE_PHYLONET/TAXON/
# MODEL_CONDITION/GENE plus the two data rows) now lands in A:E, which
# matches the target layout exactly (including which cells keep the header
# style and which shared-string values line up).
$ws.Range("A1").EntireColumn.Delete()
